$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reshape columns: drop the old column A ("INPUTS") so the remaining
# columns (old B/C widths) shift left, matching the target column widths.
$ws.Columns.Item(1).Delete()

# --- Cells that must stay TEXT even though their content looks numeric:
# mark them Text-formatted *before* writing the value so the engine keeps
# them as shared strings instead of coercing to a number.
$ws.Range("A5").NumberFormat = "@"
$ws.Range("J2").NumberFormat = "@"

# --- Header row ---
$ws.Cells.Item(1,1).Value = "TASK"
$ws.Cells.Item(1,2).Value = "ACTIVITY"
$ws.Cells.Item(1,3).Value = "CYCLES"
$ws.Cells.Item(1,4).Value = "PRIORITY"
$ws.Cells.Item(1,5).Value = "POSX"
$ws.Cells.Item(1,6).Value = "POSY"
$ws.Cells.Item(1,7).Value = ""
$ws.Cells.Item(1,8).Value = "START"
$ws.Cells.Item(1,9).Value = "CON_NAME"
$ws.Cells.Item(1,10).Value = "END"
$ws.Cells.Item(1,11).Value = ""

# --- Row 2 ---
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = "a"
$ws.Cells.Item(2,3).Value = ""
$ws.Cells.Item(2,4).Value = ""
$ws.Cells.Item(2,5).Value = ""
$ws.Cells.Item(2,6).Value = ""
$ws.Cells.Item(2,7).Value = ""
$ws.Cells.Item(2,8).Value = "1a"
$ws.Cells.Item(2,9).Value = "r1"
$ws.Cells.Item(2,10).Value = "2"
$ws.Cells.Item(2,11).Value = ""

# --- Row 3 ---
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "b"
$ws.Cells.Item(3,3).Value = ""
$ws.Cells.Item(3,4).Value = ""
$ws.Cells.Item(3,5).Value = 100
$ws.Cells.Item(3,6).Value = 150
$ws.Cells.Item(3,7).Value = ""
$ws.Cells.Item(3,8).Value = 2
$ws.Cells.Item(3,9).Value = "r2"
$ws.Cells.Item(3,10).Value = "1b"
$ws.Cells.Item(3,11).Value = ""

# --- Row 4 ---
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = ""
$ws.Cells.Item(4,3).Value = ""
$ws.Cells.Item(4,4).Value = ""
$ws.Cells.Item(4,5).Value = 300
$ws.Cells.Item(4,6).Value = 50
$ws.Cells.Item(4,7).Value = ""
$ws.Cells.Item(4,8).Value = "1a"
$ws.Cells.Item(4,9).Value = "r3"
$ws.Cells.Item(4,10).Value = "1b"
$ws.Cells.Item(4,11).Value = ""

# --- Row 5 ---
$ws.Cells.Item(5,1).Value = "3"
$ws.Cells.Item(5,2).Value = ""
$ws.Cells.Item(5,3).Value = ""
$ws.Cells.Item(5,4).Value = ""
$ws.Cells.Item(5,5).Value = ""
$ws.Cells.Item(5,6).Value = ""
$ws.Cells.Item(5,7).Value = ""
$ws.Cells.Item(5,8).Value = "1b"
$ws.Cells.Item(5,9).Value = "r4"
$ws.Cells.Item(5,10).Value = "1a"
$ws.Cells.Item(5,11).Value = ""

# --- Row 6 (only H:J populated) ---
$ws.Cells.Item(6,8).Value = 3
$ws.Cells.Item(6,9).Value = "r5"
$ws.Cells.Item(6,10).Value = 2

# --- Apply text ("@") number format -> cellXfs style 1 (numFmtId 49) to
# the whole used block. Values were written first, so numeric cells
# (counts, positions, row3/H3 etc.) keep their numeric storage while
# still picking up the new style index; A5 / J2 were pre-formatted above
# so they remain text.
$ws.Range("A1:K5").NumberFormat = "@"
$ws.Range("I6").NumberFormat = "@"

# --- Column widths ---
$ws.Columns.Item(10).ColumnWidth = 11.5

# --- Selection ---
$null = $ws.Range("K1:K6").Select()
